# Update "想去人数" (want-to-go count) values as published on gh-pages
# at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - sheet1
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 6285
$wsExpo.Range("F3").Value = 571
$wsExpo.Range("F4").Value = 124
$wsExpo.Range("F7").Value = 347
$wsExpo.Range("F8").Value = 1459

# Sheet "演出" (Performances) - sheet2
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 99

# Sheet "全部类型" (All types) - sheet4, aggregates the above sheets
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6285
$wsAll.Range("F3").Value = 571
$wsAll.Range("F4").Value = 124
$wsAll.Range("F7").Value = 347
$wsAll.Range("F8").Value = 99
$wsAll.Range("F12").Value = 1459
